$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column I (UniqueId) cells to be stored as text so the large
# 64-bit integer identifiers keep full precision (no scientific notation).
$ws.Range("I2:I71").NumberFormat = "@"

$ws.Range("I2").Value = "10690394630159179216"
$ws.Range("I3").Value = "8383772502347693165"
$ws.Range("I4").Value = "12713458539162444596"
$ws.Range("I5").Value = "17183571439025306501"
$ws.Range("I6").Value = "16291883044465602181"
$ws.Range("I7").Value = "1004739815597641664"
$ws.Range("I8").Value = "14132829364961678458"
$ws.Range("I9").Value = "15814197177705055138"
$ws.Range("I10").Value = "17153322923920234955"
$ws.Range("I11").Value = "16349991511598102228"
$ws.Range("I12").Value = "2480747059259408109"
$ws.Range("I13").Value = "15329413453490575982"
$ws.Range("I14").Value = "7877352591372251384"
$ws.Range("I15").Value = "2359784791348081432"
$ws.Range("I16").Value = "9921084047202697485"
$ws.Range("I17").Value = "15106319062049085917"
$ws.Range("I18").Value = "17887618648411249842"
$ws.Range("I19").Value = "1968137505995152539"
$ws.Range("I20").Value = "15537643808888923564"
$ws.Range("I21").Value = "8962903224833921675"
$ws.Range("I22").Value = "17839782897218791783"
$ws.Range("I23").Value = "11449381584423265720"
$ws.Range("I24").Value = "10265562418216162887"
$ws.Range("I25").Value = "5421264852785133268"
$ws.Range("I26").Value = "9775940493315330056"
$ws.Range("I27").Value = "657157722172498361"
$ws.Range("I28").Value = "13865823427859795782"
$ws.Range("I29").Value = "4032755431706838771"
$ws.Range("I30").Value = "17232621134687114681"
$ws.Range("I31").Value = "16341410120344876481"
$ws.Range("I32").Value = "5907976161846496194"
$ws.Range("I33").Value = "7975738456852880511"
$ws.Range("I34").Value = "9926735972810321660"
$ws.Range("I35").Value = "11896088144577639464"
$ws.Range("I36").Value = "3608772417383295319"
$ws.Range("I37").Value = "16143978711353994814"
$ws.Range("I38").Value = "14042058008958797846"
$ws.Range("I39").Value = "7569868334157635296"
$ws.Range("I40").Value = "9099142108530947764"
$ws.Range("I41").Value = "13654721657899965579"
$ws.Range("I42").Value = "15390026578503934594"
$ws.Range("I43").Value = "5628284592660666110"
$ws.Range("I44").Value = "7051504758948289199"
$ws.Range("I45").Value = "13907261259539446977"
$ws.Range("I46").Value = "13643525414110346129"
$ws.Range("I47").Value = "1075559904491387889"
$ws.Range("I48").Value = "4096567415949925093"
$ws.Range("I49").Value = "17630198932028886607"
$ws.Range("I50").Value = "2294922347497119376"
$ws.Range("I51").Value = "1423331922884228284"
$ws.Range("I52").Value = "10615855276681504905"
$ws.Range("I53").Value = "2138890619822242194"
$ws.Range("I54").Value = "12410837511217930970"
$ws.Range("I55").Value = "4922484447229694772"
$ws.Range("I56").Value = "9952986022782453851"
$ws.Range("I57").Value = "15752068823501477755"
$ws.Range("I58").Value = "17942581698056827285"
$ws.Range("I59").Value = "14101280709920695794"
$ws.Range("I60").Value = "12074244793306088426"
$ws.Range("I61").Value = "10198157551146233134"
$ws.Range("I62").Value = "14148964815780679069"
$ws.Range("I63").Value = "15941531297468623633"
$ws.Range("I64").Value = "2864936459970274854"
$ws.Range("I65").Value = "12329970119364176493"
$ws.Range("I66").Value = "15005498962532913587"
$ws.Range("I67").Value = "3334340575287454550"
$ws.Range("I68").Value = "6565147969932097068"
$ws.Range("I69").Value = "12047674144465933605"
$ws.Range("I70").Value = "13777818520146174391"
$ws.Range("I71").Value = "5590615921804060233"
